# Update "final projec teams.xlsx" to include all Trello board URLs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Team Game Suite -> add Trello link ---
$ws.Range("E2").Value = "https://trello.com/b/h3zbMXRF/game-development-roadmap"

# --- Row 4: Team SNES -> Trello link already present; turn it into a real hyperlink ---
$ws.Hyperlinks.Add($ws.Range("E4"), "https://trello.com/b/jVBOEz9f/team-snes-final-project")

# --- Row 5: Greenthumbs -> add Trello link ---
$ws.Range("E5").Value = "https://trello.com/b/1OGyiVu6/greenthumbs"

# --- Row 6: Team DTM -> add Trello link ---
$ws.Range("E6").Value = "https://trello.com/b/xdNUEP7R/dtm"

# --- Row 7: Dank Spots -> fix Trello link (was a stale db-setup board) ---
$ws.Range("E7").Value = "https://trello.com/b/EWzPH2WB/dankspots"

# --- Row 8: SoundBlaster -> add Angel Ruiz to team members, add Trello link ---
$ws.Range("C8").Value = "Luis Garcia, Jack Witherell, Angel Ruiz"
$ws.Range("E8").Value = "https://trello.com/b/kvv6HuKJ/cs480-music-project"

# --- Row 9: Cyber Fox Games -> add Trello link ---
$ws.Range("E9").Value = "https://trello.com/b/i3VFz7hZ/android-stealth-game"

# --- Column D (Repository) got a new explicit width once column E started filling in ---
$ws.Columns.Item(4).ColumnWidth = 18.15

# Final selection left on E10, matching the author's last-edited cell.
[void]$ws.Range("E10").Select()
